$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain stored as plain text
# (matching the source data's inlineStr cells). Force text format before
# assignment so Excel does not coerce "91.00" -> 91, "1.001" -> 1.001 (number), etc.
$textForceCells = @("D4", "D5", "D7", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D19", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = '27.766.73'
$ws.Range("E2").Value = '  -1.71%  '
$ws.Range("D3").Value = '1.893.92'
$ws.Range("E3").Value = '  -1.57%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.57%  '
$ws.Range("D5").Value = '311.90'
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("D7").Value = '0.4939'
$ws.Range("E7").Value = '  +1.75%  '
$ws.Range("E8").Value = '  -1.51%  '
$ws.Range("D9").Value = '0.07327'
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("D10").Value = '0.9096'
$ws.Range("E10").Value = '  -3.92%  '
$ws.Range("D11").Value = '20.59'
$ws.Range("E11").Value = '  -1.62%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.951.40'
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.07622'
$ws.Range("E13").Value = '  -2.23%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.459'
$ws.Range("E14").Value = '  -1.47%  '
$ws.Range("D15").Value = '6.628'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '91.00'
$ws.Range("E16").Value = '  -1.14%  '
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("D18").Value = '0.000008729'
$ws.Range("E18").Value = '  -1.72%  '
$ws.Range("D19").Value = '0.9999'
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("D20").Value = '27.648.20'
$ws.Range("E20").Value = '  -2.17%  '
$ws.Range("E21").Value = '  -3.60%  '
$ws.Range("D22").Value = '5.112'
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").Value = '2.151.76'
$ws.Range("E23").Value = '  -1.52%  '
$ws.Range("D24").Value = '10.76'
$ws.Range("E24").Value = '  -2.16%  '
$ws.Range("D25").Value = '154.03'
$ws.Range("E25").Value = '  -1.52%  '
$ws.Range("D26").Value = '1.848'
$ws.Range("E26").Value = '  -4.41%  '
$ws.Range("D27").Value = '2.183'
$ws.Range("E27").Value = '  +3.09%  '
$ws.Range("D28").Value = '18.39'
$ws.Range("E28").Value = '  -1.50%  '
$ws.Range("D29").Value = '115.07'
$ws.Range("D30").Value = '4.876'
$ws.Range("E30").Value = '  -2.87%  '
$ws.Range("D31").Value = '0.08933'
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").Value = '3.262'
$ws.Range("E32").Value = '  -3.04%  '
$ws.Range("D33").Value = '1.228'
$ws.Range("E33").Value = '  -2.11%  '
$ws.Range("D34").Value = '0.7660'
$ws.Range("E34").Value = '  -1.36%  '
$ws.Range("D35").Value = '4.640'
$ws.Range("E35").Value = '  -1.01%  '
$ws.Range("E36").Value = '  -0.42%  '
$ws.Range("D37").Value = '2.553'
$ws.Range("E37").Value = '  -7.89%  '
$ws.Range("D38").Value = '1.098'
$ws.Range("E38").Value = '  -2.98%  '
$ws.Range("E39").Value = '  -1.72%  '
$ws.Range("D40").Value = '0.05278'
$ws.Range("E40").Value = '  -1.78%  '
$ws.Range("D41").Value = '2.990'
$ws.Range("E41").Value = '  -1.79%  '
$ws.Range("D42").Value = '6.899'
$ws.Range("E42").Value = '  -2.92%  '
$ws.Range("D43").Value = '8.530'
$ws.Range("E43").Value = '  -0.72%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '0.1518'
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '112.12'
$ws.Range("E45").Value = '  +4.54%  '
$ws.Range("D46").Value = '10.58'
$ws.Range("E46").Value = '  -1.75%  '
$ws.Range("D47").Value = '0.4790'
$ws.Range("E47").Value = '  -2.73%  '
$ws.Range("D48").Value = '0.9999'
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("D49").Value = '1.629'
$ws.Range("E49").Value = '  -2.97%  '
$ws.Range("D50").Value = '67.41'
$ws.Range("E50").Value = '  -3.23%  '
$ws.Range("D51").Value = '0.06057'
$ws.Range("E51").Value = '  -1.63%  '

# Restore default (General) formatting on the forced-text cells so the
# saved XML keeps no explicit style attribute, matching the original cells.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).ClearFormats()
}
